$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") changes from 46059 -> 46060 for every data row (2-18)
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 3).Value = 46060
}

# Row 4 now holds what used to be "A 39958-2024" (previously row 17)
$ws.Cells.Item(4, 1).Value = "A 39958-2024"
$ws.Cells.Item(4, 2).Value = 45553
$ws.Cells.Item(4, 7).Value = 3.4

# Row 8 now holds what used to be "A 32256-2025" (previously row 18); loses its F (Markägare) value
$ws.Cells.Item(8, 1).Value = "A 32256-2025"
$ws.Cells.Item(8, 2).Value = 45835.6353125
$ws.Cells.Item(8, 6).Value = $null
$ws.Cells.Item(8, 7).Value = 5.4

# Row 10 now holds what used to be "A 5968-2023" (previously row 4)
$ws.Cells.Item(10, 1).Value = "A 5968-2023"
$ws.Cells.Item(10, 2).Value = 44959
$ws.Cells.Item(10, 7).Value = 1.5

# Row 11 now holds what used to be "A 2229-2023" (previously row 8)
$ws.Cells.Item(11, 1).Value = "A 2229-2023"
$ws.Cells.Item(11, 2).Value = 44939
$ws.Cells.Item(11, 7).Value = 4.3

# Row 12 now holds what used to be "A 6004-2026" (previously row 10)
$ws.Cells.Item(12, 1).Value = "A 6004-2026"
$ws.Cells.Item(12, 2).Value = 46050
$ws.Cells.Item(12, 7).Value = 2.7

# Row 13 now holds what used to be "A 5528-2023" (previously row 16)
$ws.Cells.Item(13, 1).Value = "A 5528-2023"
$ws.Cells.Item(13, 2).Value = 44957
$ws.Cells.Item(13, 7).Value = 1.2

# Row 14 now holds what used to be "A 7694-2023" (previously row 12)
$ws.Cells.Item(14, 1).Value = "A 7694-2023"
$ws.Cells.Item(14, 2).Value = 44967
$ws.Cells.Item(14, 7).Value = 2.2

# Row 16 now holds what used to be "A 39876-2024" (previously row 11)
$ws.Cells.Item(16, 1).Value = "A 39876-2024"
$ws.Cells.Item(16, 2).Value = 45553
$ws.Cells.Item(16, 7).Value = 0.3

# Row 17 now holds what used to be "A 34926-2022" (previously row 14)
$ws.Cells.Item(17, 1).Value = "A 34926-2022"
$ws.Cells.Item(17, 2).Value = 44796
$ws.Cells.Item(17, 7).Value = 1.3

# Row 18 now holds what used to be "A 28815-2024" (previously row 13); gains an F (Markägare) value
$ws.Cells.Item(18, 1).Value = "A 28815-2024"
$ws.Cells.Item(18, 2).Value = 45478
$ws.Cells.Item(18, 6).Value = "Övriga statliga verk och myndigheter"
$ws.Cells.Item(18, 7).Value = 2.8
